$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")

# Row 12
$ws.Range("H12").Value = 914.8570999999999
$ws.Range("I12").Value = 984
$ws.Range("K12").Value = 984
$ws.Range("M12").Value = -814

# Row 19
$ws.Range("H19").Value = 763.4167
$ws.Range("I19").Value = 583.6667
$ws.Range("J19").Value = 823.3333
$ws.Range("K19").Value = 583.6667
$ws.Range("L19").Value = 823.3333
$ws.Range("M19").Value = -408.6667
$ws.Range("N19").Value = -1173.3333

# Row 28
$ws.Range("H28").Value = 1014.1667
$ws.Range("I28").Value = 1072.5714
$ws.Range("J28").Value = 932.4
$ws.Range("K28").Value = 1072.5714
$ws.Range("L28").Value = 932.4
$ws.Range("M28").Value = -587.5714
$ws.Range("N28").Value = -1902.4

# Row 76
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").ClearContents()

# Row 79
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").ClearContents()

# Row 86
$ws.Range("H86").Value = 8078.7144
$ws.Range("J86").Value = 8310.4
$ws.Range("L86").Value = 8310.4
$ws.Range("N86").Value = -10556.4

# Row 89
$ws.Range("H89").Value = 8078.7144
$ws.Range("J89").Value = 8310.4
$ws.Range("L89").Value = 41552
$ws.Range("N89").Value = -52784

# Row 132
$ws.Range("H132").Value = 1088.1515
$ws.Range("I132").Value = 1106.625
$ws.Range("J132").Value = 497
$ws.Range("K132").Value = 3319.875
$ws.Range("L132").Value = 1491
$ws.Range("M132").Value = -789.875
$ws.Range("N132").Value = -6551

# Row 137
$ws.Range("H137").Value = 823.25
$ws.Range("I137").Value = 697.6667
$ws.Range("K137").Value = 2093.0001
$ws.Range("M137").Value = 456.9998999999998

# Row 138
$ws.Range("H138").Value = 3560.3914
$ws.Range("J138").Value = 4218.3125
$ws.Range("L138").Value = 12654.9375
$ws.Range("N138").Value = -22934.9375

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")

# Row 32
$ws.Range("H32").Value = 2587.56
$ws.Range("I32").Value = 2361.5
$ws.Range("J32").Value = 8013
$ws.Range("K32").Value = 2361.5
$ws.Range("L32").Value = 8013
$ws.Range("M32").Value = -2074.5
$ws.Range("N32").Value = -8587

# Row 61
$ws.Range("H61").Value = 13236.909
$ws.Range("I61").Value = 9370.571
$ws.Range("J61").Value = 20003
$ws.Range("K61").Value = 9370.571
$ws.Range("L61").Value = 20003
$ws.Range("M61").Value = -9158.571
$ws.Range("N61").Value = -20427

# Row 63
$ws.Range("H63").Value = 3316.3333
$ws.Range("I63").Value = 3339.6
$ws.Range("J63").Value = 3200
$ws.Range("K63").Value = 3339.6
$ws.Range("L63").Value = 3200
$ws.Range("M63").Value = -2653.6
$ws.Range("N63").Value = -4572

# Row 66
$ws.Range("H66").Value = 3316.3333
$ws.Range("I66").Value = 3339.6
$ws.Range("J66").Value = 3200
$ws.Range("K66").Value = 16698
$ws.Range("L66").Value = 16000
$ws.Range("M66").Value = -13266
$ws.Range("N66").Value = -22864

# Row 70
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()

# Row 73
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()

# Row 102
$ws.Range("H102").Value = 868
$ws.Range("I102").Value = 868
$ws.Range("K102").Value = 868
$ws.Range("M102").Value = 754

# Row 122
$ws.Range("H122").Value = 1348.5
$ws.Range("I122").Value = 1265.6666
$ws.Range("J122").Value = 1597
$ws.Range("K122").Value = 3796.9998
$ws.Range("L122").Value = 4791
$ws.Range("M122").Value = -1346.9998
$ws.Range("N122").Value = -9691

# Row 132
$ws.Range("H132").Value = 2846.25
$ws.Range("I132").Value = 2846.25
$ws.Range("K132").Value = 8538.75
$ws.Range("M132").Value = -6008.75

# Row 136
$ws.Range("H136").Value = 13236.909
$ws.Range("I136").Value = 9370.571
$ws.Range("J136").Value = 20003
$ws.Range("K136").Value = 28111.713
$ws.Range("L136").Value = 60009
$ws.Range("M136").Value = -25561.713
$ws.Range("N136").Value = -65109

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")

# Row 36
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()

# Row 94
$ws.Range("H94").Value = 966.5
$ws.Range("I94").Value = 966.5
$ws.Range("K94").Value = 966.5
$ws.Range("M94").Value = -515.5

# Row 107
$ws.Range("H107").Value = 699
$ws.Range("I107").Value = 699
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 699
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1221
$ws.Range("N107").ClearContents()

# Row 134
$ws.Range("H134").Value = 906.36365
$ws.Range("I134").Value = 906.36365
$ws.Range("K134").Value = 2719.09095
$ws.Range("M134").Value = -184.0909499999998

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")

# Row 31
$ws.Range("H31").Value = 1666
$ws.Range("I31").Value = 1749
$ws.Range("K31").Value = 1749
$ws.Range("M31").Value = -1454

# Row 34
$ws.Range("H34").Value = 1666
$ws.Range("I34").Value = 1749
$ws.Range("K34").Value = 1749
$ws.Range("M34").Value = -1547

# Row 58
$ws.Range("H58").Value = 3398.2856
$ws.Range("J58").Value = 2193
$ws.Range("L58").Value = 2193
$ws.Range("N58").Value = -2599

# Row 62
$ws.Range("H62").Value = 3219.8572
$ws.Range("I62").Value = 3337.8
$ws.Range("J62").Value = 2925
$ws.Range("K62").Value = 3337.8
$ws.Range("L62").Value = 2925
$ws.Range("M62").Value = -2713.8
$ws.Range("N62").Value = -4173

# Row 65
$ws.Range("H65").Value = 3219.8572
$ws.Range("I65").Value = 3337.8
$ws.Range("J65").Value = 2925
$ws.Range("K65").Value = 16689
$ws.Range("L65").Value = 14625
$ws.Range("M65").Value = -13569
$ws.Range("N65").Value = -20865

# Row 86
$ws.Range("H86").Value = 947500
$ws.Range("I86").Value = 1890000
$ws.Range("J86").Value = 5000
$ws.Range("K86").Value = 1890000
$ws.Range("L86").Value = 5000
$ws.Range("M86").Value = -1888877
$ws.Range("N86").Value = -7246

# Row 89
$ws.Range("H89").Value = 947500
$ws.Range("I89").Value = 1890000
$ws.Range("J89").Value = 5000
$ws.Range("K89").Value = 9450000
$ws.Range("L89").Value = 25000
$ws.Range("M89").Value = -9444384
$ws.Range("N89").Value = -36232

# Row 107
$ws.Range("H107").Value = 1198.5333
$ws.Range("I107").Value = 913.1539
$ws.Range("J107").Value = 3053.5
$ws.Range("K107").Value = 913.1539
$ws.Range("L107").Value = 3053.5
$ws.Range("M107").Value = 1006.8461
$ws.Range("N107").Value = -6893.5

# Row 122
$ws.Range("H122").Value = 1342
$ws.Range("I122").Value = 1342
$ws.Range("K122").Value = 4026
$ws.Range("M122").Value = -1576

# Row 136
$ws.Range("H136").Value = 3398.2856
$ws.Range("J136").Value = 2193
$ws.Range("L136").Value = 6579
$ws.Range("N136").Value = -11679

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")

# Row 131
$ws.Range("H131").Value = 1750
$ws.Range("J131").Value = 1750
$ws.Range("L131").Value = 5250
$ws.Range("N131").Value = -15330

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")

# Row 46
$ws.Range("H46").Value = 10239.833
$ws.Range("J46").Value = 12497.5
$ws.Range("L46").Value = 12497.5
$ws.Range("N46").Value = -12809.5

# Row 80
$ws.Range("H80").Value = 11003
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 11003
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 11003
$ws.Range("N80").Value = -12999
$ws.Range("M80").ClearContents()

# Row 83
$ws.Range("H83").Value = 11003
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 11003
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 55015
$ws.Range("N83").Value = -64999
$ws.Range("M83").ClearContents()

# Row 107
$ws.Range("H107").Value = 5918.6665
$ws.Range("I107").Value = 409.7143
$ws.Range("K107").Value = 409.7143
$ws.Range("M107").Value = 1510.2857

# Row 132
$ws.Range("H132").Value = 1248
$ws.Range("I132").Value = 1248
$ws.Range("K132").Value = 3744
$ws.Range("M132").Value = -1214

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")

# Row 32
$ws.Range("H32").Value = 1995
$ws.Range("I32").Value = 1995
$ws.Range("K32").Value = 1995
$ws.Range("M32").Value = -1678

# Row 61
$ws.Range("H61").Value = 2092.9
$ws.Range("I61").Value = 2104.1428
$ws.Range("K61").Value = 2104.1428
$ws.Range("M61").Value = -1902.1428

# Row 68
$ws.Range("H68").Value = 1298.6
$ws.Range("I68").Value = 1423.5
$ws.Range("J68").Value = 799
$ws.Range("K68").Value = 1423.5
$ws.Range("L68").Value = 799
$ws.Range("M68").Value = -674.5
$ws.Range("N68").Value = -2297

# Row 71
$ws.Range("H71").Value = 1298.6
$ws.Range("I71").Value = 1423.5
$ws.Range("J71").Value = 799
$ws.Range("K71").Value = 7117.5
$ws.Range("L71").Value = 3995
$ws.Range("M71").Value = -3373.5
$ws.Range("N71").Value = -11483

# Row 113
$ws.Range("H113").Value = 2092.9
$ws.Range("I113").Value = 2104.1428
$ws.Range("K113").Value = 2104.1428
$ws.Range("M113").Value = 65.85719999999992

# Row 132
$ws.Range("H132").Value = 3906.8
$ws.Range("I132").Value = 3663.1667
$ws.Range("J132").Value = 4272.25
$ws.Range("K132").Value = 10989.5001
$ws.Range("L132").Value = 12816.75
$ws.Range("M132").Value = -8459.500100000001
$ws.Range("N132").Value = -17876.75

# Row 136
$ws.Range("H136").Value = 4800
$ws.Range("I136").Value = 4800
$ws.Range("K136").Value = 14400
$ws.Range("M136").Value = -11850

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")

# Row 107
$ws.Range("H107").Value = 225.25
$ws.Range("I107").Value = 194.92308
$ws.Range("K107").Value = 584.76924
$ws.Range("M107").Value = 1335.23076

# Row 109
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

# Row 113
$ws.Range("H113").Value = 501.94446
$ws.Range("I113").Value = 573.875
$ws.Range("J113").Value = 444.4
$ws.Range("K113").Value = 1721.625
$ws.Range("L113").Value = 1333.2
$ws.Range("M113").Value = 448.375
$ws.Range("N113").Value = -5673.2

# Row 132
$ws.Range("H132").Value = 3087.5

# Row 136
$ws.Range("H136").Value = 1131.3846
$ws.Range("I136").Value = 1120.9
$ws.Range("K136").Value = 3362.7
$ws.Range("M136").Value = -812.7000000000003
